# Swap the contents of columns C (codeforiati:group-code) and D (codeforiati:group-name)
# for every row in the worksheet, including the header row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($i = 1; $i -le $rowCount; $i++) {
    $cCell = $ws.Cells.Item($i, 3)
    $dCell = $ws.Cells.Item($i, 4)
    $cVal = $cCell.Value2
    $dVal = $dCell.Value2
    $cCell.Value = $dVal
    $dCell.Value = $cVal
}
